$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 10
$ws.Range("B3").Value = 5
$ws.Range("C4").Value = 1.4

$wb.Save()
